$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Tabelle2"

# Row 12
$ws1.Range("E1").Copy()
$ws1.Range("E12:L12").PasteSpecial(-4122)
$ws1.Range("E12").Value = "A"
$ws1.Range("F12").Value = "B"
$ws1.Range("G12").Value = "C"
$ws1.Range("H12").Value = "D"
$ws1.Range("I12").Value = "E"
$ws1.Range("J12").Value = "F"
$ws1.Range("K12").Value = "G"
$ws1.Range("L12").Value = "H"

# Row 13
$ws1.Range("E1").Copy()
$ws1.Range("D13").PasteSpecial(-4122)
$ws1.Range("E2").Copy()
$ws1.Range("E13").PasteSpecial(-4122)
$ws1.Range("F2").Copy()
$ws1.Range("F13:K13").PasteSpecial(-4122)
$ws1.Range("L2").Copy()
$ws1.Range("L13").PasteSpecial(-4122)
$ws1.Range("D13").Value = 1
$ws1.Range("F13").Value = "C"
$ws1.Range("G13").Value = "A"
$ws1.Range("H13").Value = "B"
$ws1.Range("I13").Value = "B"
$ws1.Range("J13").Value = "A"
$ws1.Range("K13").Value = "C"

# Row 14
$ws1.Range("E1").Copy()
$ws1.Range("D14").PasteSpecial(-4122)
$ws1.Range("E3").Copy()
$ws1.Range("E14").PasteSpecial(-4122)
$ws1.Range("F3").Copy()
$ws1.Range("F14:K14").PasteSpecial(-4122)
$ws1.Range("L3").Copy()
$ws1.Range("L14").PasteSpecial(-4122)
$ws1.Range("D14").Value = 2
$ws1.Range("E14").Value = "C"
$ws1.Range("F14").Value = "X"
$ws1.Range("K14").Value = "X"
$ws1.Range("L14").Value = "C"

# Row 15
$ws1.Range("E1").Copy()
$ws1.Range("D15").PasteSpecial(-4122)
$ws1.Range("E3").Copy()
$ws1.Range("E15").PasteSpecial(-4122)
$ws1.Range("F3").Copy()
$ws1.Range("F15").PasteSpecial(-4122)
$ws1.Range("G4").Copy()
$ws1.Range("G15:J15").PasteSpecial(-4122)
$ws1.Range("F3").Copy()
$ws1.Range("K15").PasteSpecial(-4122)
$ws1.Range("L3").Copy()
$ws1.Range("L15").PasteSpecial(-4122)
$ws1.Range("D15").Value = 3
$ws1.Range("E15").Value = "A"
$ws1.Range("L15").Value = "A"

# Row 16
$ws1.Range("E1").Copy()
$ws1.Range("D16").PasteSpecial(-4122)
$ws1.Range("E3").Copy()
$ws1.Range("E16").PasteSpecial(-4122)
$ws1.Range("F3").Copy()
$ws1.Range("F16").PasteSpecial(-4122)
$ws1.Range("G4").Copy()
$ws1.Range("G16:J16").PasteSpecial(-4122)
$ws1.Range("F3").Copy()
$ws1.Range("K16").PasteSpecial(-4122)
$ws1.Range("L3").Copy()
$ws1.Range("L16").PasteSpecial(-4122)
$ws1.Range("D16").Value = 4
$ws1.Range("E16").Value = "B"
$ws1.Range("H16").Value = "W"
$ws1.Range("I16").Value = "S"
$ws1.Range("L16").Value = "B"

# Row 17
$ws1.Range("E1").Copy()
$ws1.Range("D17").PasteSpecial(-4122)
$ws1.Range("E3").Copy()
$ws1.Range("E17").PasteSpecial(-4122)
$ws1.Range("F3").Copy()
$ws1.Range("F17").PasteSpecial(-4122)
$ws1.Range("G4").Copy()
$ws1.Range("G17:J17").PasteSpecial(-4122)
$ws1.Range("F3").Copy()
$ws1.Range("K17").PasteSpecial(-4122)
$ws1.Range("L3").Copy()
$ws1.Range("L17").PasteSpecial(-4122)
$ws1.Range("D17").Value = 5
$ws1.Range("E17").Value = "B"
$ws1.Range("H17").Value = "S"
$ws1.Range("I17").Value = "W"
$ws1.Range("L17").Value = "B"

# Row 18
$ws1.Range("E1").Copy()
$ws1.Range("D18").PasteSpecial(-4122)
$ws1.Range("E3").Copy()
$ws1.Range("E18").PasteSpecial(-4122)
$ws1.Range("F3").Copy()
$ws1.Range("F18").PasteSpecial(-4122)
$ws1.Range("G4").Copy()
$ws1.Range("G18:J18").PasteSpecial(-4122)
$ws1.Range("F3").Copy()
$ws1.Range("K18").PasteSpecial(-4122)
$ws1.Range("L3").Copy()
$ws1.Range("L18").PasteSpecial(-4122)
$ws1.Range("D18").Value = 6
$ws1.Range("E18").Value = "A"
$ws1.Range("L18").Value = "A"

# Row 19
$ws1.Range("E1").Copy()
$ws1.Range("D19").PasteSpecial(-4122)
$ws1.Range("E3").Copy()
$ws1.Range("E19").PasteSpecial(-4122)
$ws1.Range("F3").Copy()
$ws1.Range("F19:K19").PasteSpecial(-4122)
$ws1.Range("L3").Copy()
$ws1.Range("L19").PasteSpecial(-4122)
$ws1.Range("D19").Value = 7
$ws1.Range("E19").Value = "C"
$ws1.Range("F19").Value = "X"
$ws1.Range("K19").Value = "X"
$ws1.Range("L19").Value = "C"

# Row 20
$ws1.Range("E1").Copy()
$ws1.Range("D20").PasteSpecial(-4122)
$ws1.Range("E9").Copy()
$ws1.Range("E20").PasteSpecial(-4122)
$ws1.Range("F9").Copy()
$ws1.Range("F20:K20").PasteSpecial(-4122)
$ws1.Range("L9").Copy()
$ws1.Range("L20").PasteSpecial(-4122)
$ws1.Range("D20").Value = 8
$ws1.Range("F20").Value = "C"
$ws1.Range("G20").Value = "A"
$ws1.Range("H20").Value = "B"
$ws1.Range("I20").Value = "B"
$ws1.Range("J20").Value = "A"
$ws1.Range("K20").Value = "C"

$excel.CutCopyMode = 0
$ws1.Activate()
